$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

$ws.Range("A1").Value = "number"
$ws.Range("B1").Value = "word"
$ws.Range("C1").Value = "kind"
$ws.Range("D1").Value = "carrier"
$ws.Range("E1").Value = "duplicate_image_filename"
$ws.Range("G1").Value = "order"
$ws.Range("H1").Value = "pair"
$ws.Range("I1").Value = "pair_words"
$ws.Range("J1").Value = "pair_kind"
$ws.Range("K1").Value = "carrier"

$ws.Range("A2").Value = "p1"
$ws.Range("C2").Value = "practice"
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = "A"
$ws.Range("I2").Value = "book_juice"
$ws.Range("J2").Value = "generic"
$ws.Range("K2").Value = "can"

$ws.Range("A3").Value = "p2"
$ws.Range("C3").Value = "practice"
$ws.Range("H3").Value = "B"
$ws.Range("I3").Value = "balloon_eye"
$ws.Range("J3").Value = "generic"
$ws.Range("K3").Value = "do"

$ws.Range("A4").Value = "p3"
$ws.Range("C4").Value = "practice"
$ws.Range("H4").Value = "C"
$ws.Range("I4").Value = "block_milk"
$ws.Range("J4").Value = "generic"
$ws.Range("K4").Value = "look"

$ws.Range("A5").Value = "p4"
$ws.Range("C5").Value = "practice"
$ws.Range("H5").Value = "D"
$ws.Range("I5").Value = "door_leg"
$ws.Range("J5").Value = "generic"
$ws.Range("K5").Value = "where"

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "book"
$ws.Range("C6").Value = "generic"
$ws.Range("D6").Value = "can"
$ws.Range("H6").Value = "E"
$ws.Range("K6").Value = "where"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "juice"
$ws.Range("C7").Value = "generic"
$ws.Range("D7").Value = "can"
$ws.Range("H7").Value = "F"
$ws.Range("K7").Value = "where"

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "balloon"
$ws.Range("C8").Value = "generic"
$ws.Range("D8").Value = "do"
$ws.Range("H8").Value = "G"
$ws.Range("K8").Value = "can"

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "eye"
$ws.Range("C9").Value = "generic"
$ws.Range("D9").Value = "do"
$ws.Range("H9").Value = "H"
$ws.Range("K9").Value = "do"

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "block"
$ws.Range("C10").Value = "generic"
$ws.Range("D10").Value = "look"

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "milk"
$ws.Range("C11").Value = "generic"
$ws.Range("D11").Value = "look"

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "door"
$ws.Range("C12").Value = "generic"
$ws.Range("D12").Value = "where"

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "leg"
$ws.Range("C13").Value = "generic"
$ws.Range("D13").Value = "where"

$ws.Range("A14").Value = 9

$ws.Range("A15").Value = 10

$ws.Range("A16").Value = 11

$ws.Range("A17").Value = 12

$ws.Range("A18").Value = 13

$ws.Range("A19").Value = 14

$ws.Range("A20").Value = 15

$ws.Range("A21").Value = 16

$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"

Write-Output "edit applied"
